$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 73 and 74 use the same date-column style (number format + font + border)
# as the existing last row (72); copy that cell formatting to the new date cells first.
$ws.Range("A72").Copy($ws.Range("A73"))
$ws.Range("A72").Copy($ws.Range("A74"))

# Row 73
$ws.Range("A73").Value = 45503
$ws.Range("B73").Value = 841.52151268
$ws.Range("C73").Value = 227.3684522835
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("F73").Value = 0
$ws.Range("G73").Value = 0
$ws.Range("I73").Value = 262.8307854276
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0.05384193514000001
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = 0
$ws.Range("N73").Value = 132.93283905536
$ws.Range("O73").Value = 59.76127881400001
$ws.Range("P73").Value = 0
$ws.Range("Q73").Value = [double]"2.7144E-06"
$ws.Range("R73").Value = 0
$ws.Range("S73").Value = 0
$ws.Range("T73").Value = 0
$ws.Range("U73").Value = 389.6717892303766
$ws.Range("W73").Value = 0
$ws.Range("X73").Value = 0
$ws.Range("Y73").Value = 0
$ws.Range("Z73").Value = 246.492606366138

# Row 74
$ws.Range("A74").Value = 45504
$ws.Range("B74").Value = 821.6875010800001
$ws.Range("C74").Value = 224.146392099
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("F74").Value = 0
$ws.Range("G74").Value = 0
$ws.Range("I74").Value = 251.9083454503
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0.053159309775
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = 0
$ws.Range("N74").Value = 131.39178295744
$ws.Range("O74").Value = 58.572555655
$ws.Range("P74").Value = 0
$ws.Range("Q74").Value = [double]"2.652E-06"
$ws.Range("R74").Value = 0
$ws.Range("S74").Value = 0
$ws.Range("T74").Value = 0
$ws.Range("U74").Value = 370.8662235649579
$ws.Range("W74").Value = 0
$ws.Range("X74").Value = 0
$ws.Range("Y74").Value = 0
$ws.Range("Z74").Value = 242.268921121576
